$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.23490000000001
$ws.Range("E4").Value = 16.48039999999999
$ws.Range("C7").Value = -12.53519999999999
$ws.Range("A8").Value = -22.30280000000001
$ws.Range("A10").Value = -21.47219999999998
$ws.Range("E11").Value = 16.4242
$ws.Range("A12").Value = -21.5939
$ws.Range("C14").Value = -14.39259999999999
$ws.Range("E14").Value = 16.39450000000002
$ws.Range("C15").Value = -14.36859999999998
$ws.Range("A18").Value = -21.53009999999999
$ws.Range("C18").Value = -10.3153
$ws.Range("E18").Value = 18.44450000000002
$ws.Range("E19").Value = 16.3425
$ws.Range("C20").Value = -11.7906
$ws.Range("E21").Value = 16.51760000000001
$ws.Range("A25").Value = -21.46339999999998
$ws.Range("E27").Value = 16.5112
$ws.Range("C29").Value = -11.2195
$ws.Range("C30").Value = -12.5118
$ws.Range("C31").Value = -12.96469999999999
$ws.Range("E31").Value = 16.4934
$ws.Range("C35").Value = -11.7952
$ws.Range("A37").Value = -19.57090000000001
$ws.Range("E38").Value = 16.29739999999999
$ws.Range("C40").Value = -12.6555
$ws.Range("E42").Value = 16.1879
$ws.Range("C44").Value = -12.4683
$ws.Range("E44").Value = 16.66369999999999
$ws.Range("E47").Value = 16.45360000000001
$ws.Range("C50").Value = -13.6319
$ws.Range("C54").Value = -13.1658
$ws.Range("A55").Value = -22.2684
$ws.Range("E56").Value = 16.44560000000001
$ws.Range("E58").Value = 16.14550000000002
$ws.Range("E65").Value = 17.24100000000002
$ws.Range("A68").Value = -21.66689999999999
$ws.Range("C68").Value = -11.5298
$ws.Range("E73").Value = 17.25520000000001
$ws.Range("C76").Value = -12.53060000000001
$ws.Range("A77").Value = -20.78539999999999
$ws.Range("A78").Value = -20.85789999999998
$ws.Range("A79").Value = -20.9045
$ws.Range("A80").Value = -19.4209
$ws.Range("A81").Value = -21.69000000000001
$ws.Range("A82").Value = -22.0613
$ws.Range("A84").Value = -21.967
$ws.Range("C87").Value = -13.37829999999999
$ws.Range("C88").Value = -12.75869999999999
$ws.Range("E90").Value = 16.3522
$ws.Range("C92").Value = -10.7717
$ws.Range("E92").Value = 18.29310000000002
$ws.Range("E94").Value = 18.85250000000002
$ws.Range("E95").Value = 18.10900000000002
$ws.Range("C96").Value = -12.61280000000001
$ws.Range("C98").Value = -11.94499999999999
$ws.Range("A101").Value = -21.64389999999998
$ws.Range("C101").Value = -11.50980000000001
$ws.Range("E101").Value = 16.84730000000001
$ws.Range("A102").Value = -19.27229999999999
$ws.Range("C102").Value = -13.28380000000002
